$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Determine the used range extents
$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

# Swap columns C and D (codeforiati:group-name <-> codeforiati:group-code)
# for the header row and every data row.
for ($r = 1; $r -le $lastRow; $r++) {
    $cVal = $ws.Cells.Item($r, 3).Value2
    $dVal = $ws.Cells.Item($r, 4).Value2
    $ws.Cells.Item($r, 3).Value2 = $dVal
    $ws.Cells.Item($r, 4).Value2 = $cVal
}
